$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 81;  "F2" = 54;  "H2" = 66;
    "E4" = 32;
    "E9" = 38;
    "E12" = 14; "F12" = 8;  "H12" = 8;
    "E15" = 200;
    "E17" = 155; "F17" = 82; "H17" = 114;
    "E18" = 152; "F18" = 74; "H18" = 111;
    "E19" = 87;
    "E23" = 12;
    "E24" = 32;
    "E25" = 32; "F25" = 21; "H25" = 29;
    "E26" = 44;
    "E34" = 33;
    "F35" = 10; "H35" = 11;
    "E37" = 75; "F37" = 45; "H37" = 57;
    "E38" = 103; "F38" = 29; "H38" = 49;
    "E41" = 54; "F41" = 29; "H41" = 40;
    "E43" = 40;
    "E47" = 79; "F47" = 52; "H47" = 62;
    "E49" = 83;
    "E51" = 16; "F51" = 12; "H51" = 12;
    "E55" = 11;
    "E56" = 12; "F56" = 5;  "H56" = 7;
    "E57" = 26; "F57" = 7;  "H57" = 11;
    "E59" = 13; "F59" = 7;  "H59" = 11;
    "E60" = 26;
    "F61" = 19; "H61" = 29;
    "E62" = 66;
    "E64" = 47; "F64" = 24; "H64" = 29;
    "E65" = 45; "F65" = 18; "H65" = 31;
    "E66" = 44;
    "E67" = 51;
    "E68" = 20;
    "E69" = 22;
    "E71" = 55;
    "E72" = 58; "F72" = 32; "H72" = 43;
    "E73" = 42; "F73" = 22; "H73" = 34;
    "F76" = 25; "H76" = 42;
    "E77" = 70; "F77" = 30; "H77" = 47;
    "F78" = 25; "H78" = 46;
    "E79" = 53; "F79" = 28; "H79" = 39;
    "E83" = 15;
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
